# "Improve BLE characteristic details"
#
# Updates the "BLE characteristics" sheet (sheet1) with corrected /
# expanded descriptions for the dc_serv characteristics, and adds new
# rows describing the burst_serv characteristics (pulses, bursts,
# p1_curr, p2_curr, anodic, timing). Also swaps which sheet/cell is the
# active selection, and widens column D to fit the longer text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- dc_serv rows: reworded descriptions (row 5-7) ---
$ws1.Range("D5").Value = "4 byte array. Float referring to target value of DC stim, current. Little endian format."
$ws1.Range("D6").Value = "4 byte array. Float referring to base value of DC stim, current. Little endian format."
$ws1.Range("D7").Value = "12 byte array. 3 sets of 4 bytes, each in little endian format. Order: Slope, Flat, Interstim"

# --- finish row 8 (pulses) which previously only had B8 set ---
$ws1.Range("C8").Value = "read"
$ws1.Range("D8").Value = "4 byte array. Little endian uint32_t. Number of pulses we are aiming to complete with DC"

# --- row 9: burst_serv service header + first characteristic (pulses) ---
$ws1.Range("A9").Value = "burst_serv"
$ws1.Range("B9").Value = "pulses"
$ws1.Range("C9").Value = "read"
$ws1.Range("D9").Value = "4 byte array. Little endian uint32_t"

# --- row 10: bursts ---
$ws1.Range("B10").Value = "bursts"
$ws1.Range("C10").Value = "read"
$ws1.Range("D10").Value = "4 byte array. Little endian uint32_t"

# --- row 11: p1_curr ---
$ws1.Range("B11").Value = "p1_curr"
$ws1.Range("C11").Value = "read"
$ws1.Range("D11").Value = "4 byte array. Little endian float giving current of phase 1"

# --- row 12: p2_curr ---
$ws1.Range("B12").Value = "p2_curr"
$ws1.Range("C12").Value = "read"
$ws1.Range("D12").Value = "4 byte array. Little endian float giving current of phase 2"

# --- row 13: anodic ---
$ws1.Range("B13").Value = "anodic"
$ws1.Range("C13").Value = "read"
$ws1.Range("D13").Value = "1 byte array. 1 if anodic, 0 if not."

# --- row 14: timing (burst_serv variant) ---
$ws1.Range("B14").Value = "timing"
$ws1.Range("C14").Value = "read"
$ws1.Range("D14").Value = "20 byte array. 5 sets of 4 bytes, each set in little endian format. Order: Interstim, P1, Interphase, P2, Interburst"

# --- widen column D to fit the new, longer text ---
$ws1.Columns.Item(4).ColumnWidth = 99.6

# --- swap the active sheet / selection: "BLE characteristics" becomes
#     the active tab (selecting D5), "Command details" no longer is
#     (selecting D5 there too) ---
$ws2.Range("D5").Select()
$ws1.Activate()
$ws1.Range("D5").Select()
